$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07702833333333334
$ws.Range("H2").Value = 0.231085
$ws.Range("I2").Value = 0.02259036512642383
$ws.Range("J2").Value = 0.02259036512642383
$ws.Range("O2").Value = 0.7140239834365498
$ws.Range("P2").Value = 0.7140239834365498
$ws.Range("Q2").Value = 0.1542426387394444
$ws.Range("R2").Value = 1.388183748655
$ws.Range("S2").Value = 0.01613006249485526
$ws.Range("T2").Value = 0.01613006249485526

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07702833333333334
$ws.Range("H3").Value = 0.231085
$ws.Range("I3").Value = 0.02259036512642383
$ws.Range("J3").Value = 0.02259036512642383
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8019933333333333
$ws.Range("N3").Value = 2.40598
$ws.Range("O3").Value = 0.2859760165634502
$ws.Range("P3").Value = 0.2859760165634502
$ws.Range("Q3").Value = 0.06177620981111111
$ws.Range("R3").Value = 0.5559858883000001
$ws.Range("S3").Value = 0.006460302631568569
$ws.Range("T3").Value = 0.006460302631568569

$ws.Range("I4").Value = 0.9608869019286738
$ws.Range("J4").Value = 0.9608869019286738
$ws.Range("O4").Value = 0.7140239834365498
$ws.Range("P4").Value = 0.7140239834365498
$ws.Range("S4").Value = 0.6860962933471171
$ws.Range("T4").Value = 0.6860962933471171

$ws.Range("I5").Value = 0.9608869019286738
$ws.Range("J5").Value = 0.9608869019286738
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8019933333333333
$ws.Range("N5").Value = 2.40598
$ws.Range("O5").Value = 0.2859760165634502
$ws.Range("P5").Value = 0.2859760165634502
$ws.Range("Q5").Value = 2.627666729868889
$ws.Range("R5").Value = 23.64900056882
$ws.Range("S5").Value = 0.2747906085815567
$ws.Range("T5").Value = 0.2747906085815567

$ws.Range("G6").Value = 0.056339
$ws.Range("H6").Value = 0.169017
$ws.Range("I6").Value = 0.01652273294490242
$ws.Range("J6").Value = 0.01652273294490242
$ws.Range("O6").Value = 0.7140239834365498
$ws.Range("P6").Value = 0.7140239834365498
$ws.Range("Q6").Value = 0.1128140211256667
$ws.Range("R6").Value = 1.015326190131
$ws.Range("S6").Value = 0.01179762759457754
$ws.Range("T6").Value = 0.01179762759457754

$ws.Range("G7").Value = 0.056339
$ws.Range("H7").Value = 0.169017
$ws.Range("I7").Value = 0.01652273294490242
$ws.Range("J7").Value = 0.01652273294490242
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8019933333333333
$ws.Range("N7").Value = 2.40598
$ws.Range("O7").Value = 0.2859760165634502
$ws.Range("P7").Value = 0.2859760165634502
$ws.Range("Q7").Value = 0.04518350240666667
$ws.Range("R7").Value = 0.40665152166
$ws.Range("S7").Value = 0.004725105350324879
$ws.Range("T7").Value = 0.004725105350324879
